$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 255, shifting existing rows 255:271 down to 256:272
$ws.Rows.Item(255).Insert()

# The newly inserted row 255 inherited formatting from the row above (254),
# but we need it to match the row that is now below it (256, formerly 255).
# Copy the previous row 256 (post-shift, originally row 255) into the new row 255
# so all static columns/styles line up, then overwrite the changed values.
$ws.Rows.Item(256).Copy()
$ws.Rows.Item(255).PasteSpecial()
$excel.CutCopyMode = $false

# Now set the new data values for row 255
$ws.Range("D255").Value = 44706
$ws.Range("E255").Value = 4
$ws.Range("F255").Value = 100112012
$ws.Range("G255").Value = "Espinaca"
$ws.Range("H255").Value = "Sin especificar"
$ws.Range("I255").Value = "Primera"
$ws.Range("J255").Value = 2800
$ws.Range("K255").Value = 500
$ws.Range("L255").Value = 600
$ws.Range("M255").Value = 550
$ws.Range("N255").Value = "$/atado 300 a 500 gramos"
$ws.Range("O255").Value = "Provincia del Elquí"
$ws.Range("P255").Value = 1100
$ws.Range("Q255").Value = 0.5
$ws.Range("R255").Value = "Hortaliza"
$ws.Range("A255").Value = 8
$ws.Range("B255").Value = "Terminal La Palmera de La Serena"
$ws.Range("C255").Value = "Coquimbo"
